# Update the cryptos list with latest scraped values (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($rangeAddr, $text) {
    # Force the cell to stay text even when the new value looks like a
    # plain number (e.g. "216.41"), matching the source data which stores
    # all Price/Volume column entries as text. Reset the style back to
    # "Normal" afterwards so no residual number-format style lingers on
    # the cell (keeps the style index identical to the untouched cells).
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.042.15"
$ws.Range("E2").Value = "  +0.40%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.644.10"
$ws.Range("E3").Value = "  +0.67%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.85%  "

# Row 5 - BNB
Set-TextValue "D5" "216.41"
$ws.Range("E5").Value = "  +0.70%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.84%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.74%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +0.45%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +1.11%  "

# Row 10 - Solana
Set-TextValue "D10" "19.68"
$ws.Range("E10").Value = "  -0.03%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.98%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.872.44"
$ws.Range("E12").Value = "  +0.77%  "

# Row 13 - Polkadot
$ws.Range("E13").Value = "  +1.32%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.650.76"
$ws.Range("E14").Value = "  +1.83%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -0.03%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +0.79%  "

# Row 17 - Litecoin
$ws.Range("E17").Value = "  +0.58%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "26.041.79"
$ws.Range("E18").Value = "  +0.57%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.78%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "193.08"
$ws.Range("E20").Value = "  +0.04%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -0.85%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  -0.07%  "

# Row 23 - Chainlink
Set-TextValue "D23" "6.24"
$ws.Range("E23").Value = "  -0.32%  "

# Row 24 - now Toncoin (was Stellar)
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D24" "1.80"
$ws.Range("E24").Value = "  +0.75%  "

# Row 25 - now Stellar (was Monero)
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D25" "0.131"
$ws.Range("E25").Value = "  +4.24%  "

# Row 26 - now Monero (was Toncoin)
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D26" "144.59"
$ws.Range("E26").Value = "  +1.52%  "

# Row 27 - BinanceUSD
$ws.Range("E27").Value = "  +0.95%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  +0.35%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "15.53"
$ws.Range("E29").Value = "  +0.44%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.98%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -0.19%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  -0.63%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  +1.09%  "

# Row 34 - now HuobiToken (was LidoDAOToken)
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D34" "2.47"
$ws.Range("E34").Value = "  +2.43%  "

# Row 35 - now LidoDAOToken (was HuobiToken)
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D35" "1.53"
$ws.Range("E35").Value = "  -3.22%  "

# Row 36 - ARBITRUM
$ws.Range("E36").Value = "  +0.38%  "

# Row 37 - Maker
$ws.Range("D37").Value = "1.132.88"
$ws.Range("E37").Value = "  -0.33%  "

# Row 38 - ImmutableX
Set-TextValue "D38" "0.540"
$ws.Range("E38").Value = "  -1.75%  "

# Row 39 - MXToken
$ws.Range("E39").Value = "  +0.45%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +0.46%  "

# Row 41 - FraxShare
$ws.Range("E41").Value = "  +0.65%  "

# Row 42 - Quant
Set-TextValue "D42" "99.46"
$ws.Range("E42").Value = "  +0.29%  "

# Row 43 - TrustWalletToken
Set-TextValue "D43" "0.797"
$ws.Range("E43").Value = "  -0.74%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.782.03"
$ws.Range("E44").Value = "  +0.87%  "

# Row 45 - BabyDogeCoin
$ws.Range("D45").Value = "0.0₆0115"
$ws.Range("E45").Value = "  +3.42%  "

# Row 46 - Aave
$ws.Range("E46").Value = "  +0.69%  "

# Row 47 - Cronos
Set-TextValue "D47" "0.0531"
$ws.Range("E47").Value = "  -0.06%  "

# Row 48 - RenderToken
Set-TextValue "D48" "1.45"
$ws.Range("E48").Value = "  -0.21%  "

# Row 49 - EnergySwap
Set-TextValue "D49" "7.69"
$ws.Range("E49").Value = "  +0.77%  "

# Row 50 - Mantle
Set-TextValue "D50" "0.417"
$ws.Range("E50").Value = "  +0.57%  "

# Row 51 - Algorand
$ws.Range("E51").Value = "  -0.53%  "
